$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing text in A46: "14/05/2020" -> "14/5/2020"
$ws.Range("A46").Value = "14/5/2020"

# Add new row 47 data
$ws.Range("A47").Value = "15/5/2020"
$ws.Range("B47").Value = 1233
$ws.Range("C47").Value = 125

# Match style of the row above (right-aligned "General" style) for the date text cell
$ws.Range("A47").HorizontalAlignment = -4152  # xlRight

# Update the selection to match the diff (active cell C48)
$ws.Range("C48").Select()

# Resize the table to include the new row
$tbl = $ws.ListObjects.Item("Condicion_Pacientes")
$tbl.Resize($ws.Range("A1:F47"))
